$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) from "ShearF-HW25.xpc" to "ShearF"
$ws.Name = "ShearF"

# Add new row 16 with averaged intensity data for HKL entry 14,
# using the same style (s=1) as column A in the other rows (copy the
# formatting from A15 so the existing style record is reused rather
# than a brand-new one created) and the "HexGrid-60degTilt5degRes"
# label (shared string) for column B.
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 0.9519747158388456
$ws.Cells.Item(16, 4).Value = 1.342963186989
$ws.Cells.Item(16, 5).Value = 0.8977606274731419
$ws.Cells.Item(16, 6).Value = 0.9519747158388456
$ws.Cells.Item(16, 7).Value = 1.176649270150119
$ws.Cells.Item(16, 8).Value = 0.7618749660278448
$ws.Cells.Item(16, 9).Value = 0.9090833387038024
$ws.Cells.Item(16, 10).Value = 1.342963186989
$ws.Cells.Item(16, 11).Value = 1.120361907231071
$ws.Cells.Item(16, 12).Value = 1.036168311534958
$ws.Cells.Item(16, 13).Value = 1.006717684197125
